$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 710.8889
$ws.Range("I15").Value = 710.8889
$ws.Range("K15").Value = 2132.6667
$ws.Range("M15").Value = -1963.6667
$ws.Range("H17").Value = 630.9143
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 630.9143
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 1892.7429
$ws.Range("M17").Value = $null
$ws.Range("N17").Value = -2228.7429
$ws.Range("H43").Value = 4956.558
$ws.Range("J43").Value = 4616.9546
$ws.Range("L43").Value = 4616.9546
$ws.Range("N43").Value = -4754.9546
$ws.Range("H64").Value = 9548.5
$ws.Range("J64").Value = 9548.5
$ws.Range("L64").Value = 9548.5
$ws.Range("N64").Value = -10044.5
$ws.Range("H67").Value = 9548.5
$ws.Range("J67").Value = 9548.5
$ws.Range("L67").Value = 9548.5
$ws.Range("N67").Value = -11264.5
$ws.Range("H70").Value = 2997
$ws.Range("I70").Value = 2800
$ws.Range("K70").Value = 8400
$ws.Range("M70").Value = -8130
$ws.Range("H73").Value = 2997
$ws.Range("I73").Value = 2800
$ws.Range("K73").Value = 8400
$ws.Range("M73").Value = -7464
$ws.Range("H99").Value = 2528.2307
$ws.Range("J99").Value = 4114.6665
$ws.Range("L99").Value = 12343.9995
$ws.Range("N99").Value = -15339.9995
$ws.Range("H101").Value = 1789.5714
$ws.Range("J101").Value = 4095
$ws.Range("L101").Value = 12285
$ws.Range("N101").Value = -15529
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").Value = $null
$ws.Range("H127").Value = 588.125
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").Value = $null
$ws.Range("H131").Value = 5747.5
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").Value = $null
$ws.Range("H132").Value = 47107.184
$ws.Range("I132").Value = 51677.35
$ws.Range("K132").Value = 155032.05
$ws.Range("M132").Value = -152502.05
$ws.Range("H138").Value = 7048.544
$ws.Range("I138").Value = 3762.1538
$ws.Range("J138").Value = 8019.523
$ws.Range("K138").Value = 11286.4614
$ws.Range("L138").Value = 24058.569
$ws.Range("M138").Value = -6146.4614
$ws.Range("N138").Value = -34338.569

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1844.238
$ws.Range("I45").Value = 1621.3125
$ws.Range("J45").Value = 2557.6
$ws.Range("K45").Value = 1621.3125
$ws.Range("L45").Value = 2557.6
$ws.Range("M45").Value = -1244.3125
$ws.Range("N45").Value = -3311.6
$ws.Range("H61").Value = 1896.8334
$ws.Range("I61").Value = 1997.75
$ws.Range("K61").Value = 1997.75
$ws.Range("M61").Value = -1785.75
$ws.Range("H136").Value = 1896.8334
$ws.Range("I136").Value = 1997.75
$ws.Range("K136").Value = 5993.25
$ws.Range("M136").Value = -3443.25

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 429.77777
$ws.Range("I22").Value = 397.66666
$ws.Range("K22").Value = 397.66666
$ws.Range("M22").Value = -224.66666

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2914.3333
$ws.Range("I16").Value = 2778.625
$ws.Range("K16").Value = 2778.625
$ws.Range("M16").Value = -2491.625
$ws.Range("H31").Value = 2577.9524
$ws.Range("J31").Value = 2222.0588
$ws.Range("L31").Value = 2222.0588
$ws.Range("N31").Value = -2812.0588
$ws.Range("H34").Value = 2577.9524
$ws.Range("J34").Value = 2222.0588
$ws.Range("L34").Value = 2222.0588
$ws.Range("N34").Value = -2626.0588
$ws.Range("H94").Value = 729
$ws.Range("I94").Value = 484.8
$ws.Range("J94").Value = 903.4286
$ws.Range("K94").Value = 484.8
$ws.Range("L94").Value = 903.4286
$ws.Range("M94").Value = -33.80000000000001
$ws.Range("N94").Value = -1805.4286
$ws.Range("H107").Value = 5294.6665
$ws.Range("I107").Value = 2351.6155
$ws.Range("K107").Value = 2351.6155
$ws.Range("M107").Value = -431.6154999999999
$ws.Range("H113").Value = 2914.3333
$ws.Range("I113").Value = 2778.625
$ws.Range("K113").Value = 2778.625
$ws.Range("M113").Value = -608.625
$ws.Range("H132").Value = 3142.15
$ws.Range("I132").Value = 3216.6428
$ws.Range("J132").Value = 2968.3333
$ws.Range("K132").Value = 9649.928400000001
$ws.Range("L132").Value = 8904.999899999999
$ws.Range("M132").Value = -7119.928400000001
$ws.Range("N132").Value = -13964.9999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 1641.7693
$ws.Range("I2").Value = 2325.1667
$ws.Range("K2").Value = 13951.0002
$ws.Range("M2").Value = -13838.0002
$ws.Range("H5").Value = 1939.2
$ws.Range("I5").Value = 1232.8334
$ws.Range("K5").Value = 3698.5002
$ws.Range("M5").Value = -3586.5002
$ws.Range("H26").Value = 612.0526
$ws.Range("I26").Value = 287.5
$ws.Range("J26").Value = 698.6
$ws.Range("K26").Value = 862.5
$ws.Range("L26").Value = 2095.8
$ws.Range("M26").Value = -574.5
$ws.Range("N26").Value = -2671.8
$ws.Range("H47").Value = 3289
$ws.Range("I47").Value = 1050.25
$ws.Range("J47").Value = 5080
$ws.Range("K47").Value = 3150.75
$ws.Range("L47").Value = 15240
$ws.Range("M47").Value = -2719.75
$ws.Range("N47").Value = -16102
$ws.Range("H68").Value = 1257.6316
$ws.Range("I68").Value = 920.55554
$ws.Range("J68").Value = 1561
$ws.Range("K68").Value = 2761.66662
$ws.Range("L68").Value = 4683
$ws.Range("M68").Value = -1950.66662
$ws.Range("N68").Value = -6305
$ws.Range("H71").Value = 1257.6316
$ws.Range("I71").Value = 920.55554
$ws.Range("J71").Value = 1561
$ws.Range("K71").Value = 8284.99986
$ws.Range("L71").Value = 14049
$ws.Range("M71").Value = -4228.99986
$ws.Range("N71").Value = -22161
$ws.Range("H107").Value = 1730.909
$ws.Range("J107").Value = 2049.611
$ws.Range("L107").Value = 6148.833
$ws.Range("N107").Value = -9988.832999999999
$ws.Range("H135").Value = 1939.2
$ws.Range("I135").Value = 1232.8334
$ws.Range("K135").Value = 11095.5006
$ws.Range("M135").Value = -8560.500599999999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2166.3125
$ws.Range("I102").Value = 1789.25
$ws.Range("K102").Value = 1789.25
$ws.Range("M102").Value = -167.25
$ws.Range("H107").Value = 928.2105
$ws.Range("I107").Value = 926.2941
$ws.Range("J107").Value = 944.5
$ws.Range("K107").Value = 926.2941
$ws.Range("L107").Value = 944.5
$ws.Range("M107").Value = 993.7059
$ws.Range("N107").Value = -4784.5
$ws.Range("H126").Value = 4201.5386
$ws.Range("I126").Value = 3404.4443
$ws.Range("J126").Value = 5995
$ws.Range("K126").Value = 10213.3329
$ws.Range("L126").Value = 17985
$ws.Range("M126").Value = -7743.332900000001
$ws.Range("N126").Value = -22925

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").Value = $null
$ws.Range("H124").Value = 54724.5
$ws.Range("J124").Value = 54724.5
$ws.Range("L124").Value = 54724.5
$ws.Range("N124").Value = -64544.5
$ws.Range("H132").Value = 3167.7778
$ws.Range("I132").Value = 3001.4285
$ws.Range("J132").Value = 3750
$ws.Range("K132").Value = 9004.2855
$ws.Range("L132").Value = 11250
$ws.Range("M132").Value = -6474.2855
$ws.Range("N132").Value = -16310

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 11899
$ws.Range("I14").Value = 6990
$ws.Range("J14").Value = 13126.25
$ws.Range("K14").Value = 6990
$ws.Range("L14").Value = 13126.25
$ws.Range("M14").Value = -6822
$ws.Range("N14").Value = -13462.25
$ws.Range("H27").Value = 24198.5
$ws.Range("J27").Value = 24198.5
$ws.Range("L27").Value = 24198.5
$ws.Range("N27").Value = -24336.5
$ws.Range("H62").Value = 3099.75
$ws.Range("I62").Value = 1949.5
$ws.Range("J62").Value = 4250
$ws.Range("K62").Value = 1949.5
$ws.Range("L62").Value = 4250
$ws.Range("M62").Value = -1325.5
$ws.Range("N62").Value = -5498
$ws.Range("H65").Value = 3099.75
$ws.Range("I65").Value = 1949.5
$ws.Range("J65").Value = 4250
$ws.Range("K65").Value = 9747.5
$ws.Range("L65").Value = 21250
$ws.Range("M65").Value = -6627.5
$ws.Range("N65").Value = -27490
$ws.Range("H122").Value = 2067.158
$ws.Range("J122").Value = 3005
$ws.Range("L122").Value = 9015
$ws.Range("N122").Value = -13915
$ws.Range("H132").Value = 3867.182
$ws.Range("I132").Value = 3675.2942
$ws.Range("K132").Value = 11025.8826
$ws.Range("M132").Value = -8495.882599999999
